# Auto-generated Excel COM script to refresh Famfrit_Profits market-data columns (H-N)
# across all 8 item-sourcing sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 1999
$ws.Range("I26").Value = 1999
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1999
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -1655
$ws.Range("N26").ClearContents()
# Row 43
$ws.Range("H43").Value = 10000000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
# Row 48
$ws.Range("H48").Value = 3281.3333
$ws.Range("J48").Value = 3172
$ws.Range("L48").Value = 9516
$ws.Range("N48").Value = -10100
# Row 56
$ws.Range("H56").Value = 3281.3333
$ws.Range("J56").Value = 3172
$ws.Range("L56").Value = 9516
$ws.Range("N56").Value = -10584
# Row 76
$ws.Range("H76").Value = 12125.77
$ws.Range("J76").Value = 4678.5713
$ws.Range("L76").Value = 4678.5713
$ws.Range("N76").Value = -5308.5713
# Row 79
$ws.Range("H79").Value = 12125.77
$ws.Range("J79").Value = 4678.5713
$ws.Range("L79").Value = 4678.5713
$ws.Range("N79").Value = -6862.5713
# Row 80
$ws.Range("H80").Value = 2245.3103
$ws.Range("I80").Value = 1582.4166
$ws.Range("K80").Value = 4747.2498
$ws.Range("M80").Value = -3749.2498
# Row 83
$ws.Range("H83").Value = 2245.3103
$ws.Range("I83").Value = 1582.4166
$ws.Range("K83").Value = 14241.7494
$ws.Range("M83").Value = -9249.749400000001
# Row 98
$ws.Range("H98").Value = 1140.7778
$ws.Range("I98").Value = 1064.625
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 1064.625
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = 433.375
$ws.Range("N98").Value = -4746
# Row 122
$ws.Range("H122").Value = 1140.7778
$ws.Range("I122").Value = 1064.625
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3193.875
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -743.875
$ws.Range("N122").Value = -10150
# Row 132
$ws.Range("H132").Value = 2324.742
$ws.Range("I132").Value = 2070.8965
$ws.Range("K132").Value = 6212.689499999999
$ws.Range("M132").Value = -3682.689499999999
# Row 138
$ws.Range("H138").Value = 9010510
$ws.Range("I138").Value = 767.25
$ws.Range("J138").Value = 19610208
$ws.Range("K138").Value = 2301.75
$ws.Range("L138").Value = 58830624
$ws.Range("M138").Value = 2838.25
$ws.Range("N138").Value = -58840904

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 5191.769
$ws.Range("I45").Value = 4275.5557
$ws.Range("K45").Value = 4275.5557
$ws.Range("M45").Value = -3898.5557
# Row 74
$ws.Range("H74").Value = 71430590
$ws.Range("I74").Value = 142858700
$ws.Range("K74").Value = 142858700
$ws.Range("M74").Value = -142857826
# Row 77
$ws.Range("H77").Value = 71430590
$ws.Range("I77").Value = 142858700
$ws.Range("K77").Value = 714293500
$ws.Range("M77").Value = -714289132
# Row 122
$ws.Range("H122").Value = 2717.1292
$ws.Range("I122").Value = 1892.35
$ws.Range("K122").Value = 5677.049999999999
$ws.Range("M122").Value = -3227.049999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 7002.7646
$ws.Range("I105").Value = 9187.333000000001
$ws.Range("K105").Value = 9187.333000000001
$ws.Range("M105").Value = -7440.333000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7365.273
$ws.Range("I31").Value = 3750.8333
$ws.Range("K31").Value = 3750.8333
$ws.Range("M31").Value = -3455.8333
# Row 34
$ws.Range("H34").Value = 7365.273
$ws.Range("I34").Value = 3750.8333
$ws.Range("K34").Value = 3750.8333
$ws.Range("M34").Value = -3548.8333
# Row 134
$ws.Range("H134").Value = 2340
$ws.Range("I134").Value = 2233.3333
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 6699.999899999999
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -4164.999899999999
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 73
$ws.Range("J11").Value = 80
$ws.Range("L11").Value = 240
$ws.Range("N11").Value = -520
# Row 75
$ws.Range("H75").Value = 1979.8
$ws.Range("I75").Value = 1999
$ws.Range("J75").Value = 1975
$ws.Range("K75").Value = 5997
$ws.Range("L75").Value = 5925
$ws.Range("M75").Value = -4999
$ws.Range("N75").Value = -7921
# Row 78
$ws.Range("H78").Value = 1979.8
$ws.Range("I78").Value = 1999
$ws.Range("J78").Value = 1975
$ws.Range("K78").Value = 17991
$ws.Range("L78").Value = 17775
$ws.Range("M78").Value = -12999
$ws.Range("N78").Value = -27759
# Row 128
$ws.Range("H128").Value = 120995.664
$ws.Range("I128").Value = 120995.664
$ws.Range("K128").Value = 362986.992
$ws.Range("M128").Value = -358006.992
# Row 131
$ws.Range("H131").Value = 30325.666
$ws.Range("I131").Value = 116224.89
$ws.Range("J131").Value = 4555.9
$ws.Range("K131").Value = 348674.67
$ws.Range("L131").Value = 13667.7
$ws.Range("M131").Value = -343634.67
$ws.Range("N131").Value = -23747.7

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 114739.39
$ws.Range("I70").Value = 121253.53
$ws.Range("J70").Value = 3999
$ws.Range("K70").Value = 121253.53
$ws.Range("L70").Value = 3999
$ws.Range("M70").Value = -120983.53
$ws.Range("N70").Value = -4539
# Row 73
$ws.Range("H73").Value = 114739.39
$ws.Range("I73").Value = 121253.53
$ws.Range("J73").Value = 3999
$ws.Range("K73").Value = 121253.53
$ws.Range("L73").Value = 3999
$ws.Range("M73").Value = -120317.53
$ws.Range("N73").Value = -5871
# Row 97
$ws.Range("H97").Value = 1916.4615
$ws.Range("J97").Value = 4627.75
$ws.Range("L97").Value = 4627.75
$ws.Range("N97").Value = -5619.75
# Row 122
$ws.Range("H122").Value = 1802.6
$ws.Range("I122").Value = 1561.5834
$ws.Range("K122").Value = 4684.7502
$ws.Range("M122").Value = -2234.7502
# Row 126
$ws.Range("H126").Value = 3295.4546
$ws.Range("I126").Value = 1382
$ws.Range("K126").Value = 4146
$ws.Range("M126").Value = -1676

$ws = $wb.Worksheets.Item("LTW")
# Row 30
$ws.Range("H30").Value = 2275
$ws.Range("I30").Value = 2275
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2275
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2167
$ws.Range("N30").ClearContents()
# Row 35
$ws.Range("H35").Value = 650
$ws.Range("I35").Value = 650
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 650
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -314
$ws.Range("N35").ClearContents()
# Row 40
$ws.Range("H40").Value = 2987.7222
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 82
$ws.Range("H82").Value = 2559.484
$ws.Range("I82").Value = 1561
$ws.Range("K82").Value = 1561
$ws.Range("M82").Value = -1200
# Row 85
$ws.Range("H85").Value = 2559.484
$ws.Range("I85").Value = 1561
$ws.Range("K85").Value = 1561
$ws.Range("M85").Value = -313
# Row 122
$ws.Range("H122").Value = 4773.41
$ws.Range("I122").Value = 3891.3572
$ws.Range("K122").Value = 11674.0716
$ws.Range("M122").Value = -9224.071599999999
# Row 132
$ws.Range("H132").Value = 10081.85
$ws.Range("I132").Value = 7664.136
$ws.Range("K132").Value = 22992.408
$ws.Range("M132").Value = -20462.408

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1550.75
$ws.Range("I81").Value = 1367.3334
$ws.Range("K81").Value = 2734.6668
$ws.Range("M81").Value = -1673.6668
# Row 84
$ws.Range("H84").Value = 1550.75
$ws.Range("I84").Value = 1367.3334
$ws.Range("K84").Value = 13673.334
$ws.Range("M84").Value = -8369.333999999999
# Row 122
$ws.Range("H122").Value = 69423.734
$ws.Range("J122").Value = 4246.75
$ws.Range("L122").Value = 12740.25
$ws.Range("N122").Value = -17640.25
